# Reorders the weekly price records in rows 3-16 (the header is row 1 and
# row 2 is untouched). Columns A, B, C, E, F, G, H, I, N, Q, R are identical
# across every data row, so only the per-record columns D (Fecha), J
# (Volumen), K (Precio minimo), L (Precio maximo), M (Precio promedio
# ponderado), O (Origen) and P (Precio $/Kg) need to move between rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry per-record data which gets reshuffled.
$cols = @(4, 10, 11, 12, 13, 15, 16)   # D, J, K, L, M, O, P

# Maps each destination row (after the edit) to the row whose data it
# should receive (before the edit), for rows 3..16.
$rowMap = @{
  3  = 5
  4  = 10
  5  = 16
  6  = 15
  7  = 8
  8  = 4
  9  = 13
  10 = 3
  11 = 14
  12 = 11
  13 = 9
  14 = 12
  15 = 7
  16 = 6
}

# Snapshot the original values for every row/column involved before writing
# anything, since several rows trade values with each other.
$snapshot = @{}
foreach ($r in $rowMap.Keys) {
  $rowVals = @{}
  foreach ($c in $cols) {
    $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
  }
  $snapshot[$r] = $rowVals
}

foreach ($destRow in $rowMap.Keys) {
  $srcRow = $rowMap[$destRow]
  $srcVals = $snapshot[$srcRow]
  foreach ($c in $cols) {
    $ws.Cells.Item($destRow, $c).Value2 = $srcVals[$c]
  }
}
